$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 (introduces "estrella")
$ws.Range("G1").Value = "estrella"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# New row 12
$ws.Range("A12").Value = "autSviRep"
$ws.Range("B12").Value = "Autorización SVI representante"
$ws.Range("C12").Value = "Autorización SVI representante"
$ws.Range("G12").Value = "DEN_SVI_REP"
$ws.Range("F12").Value = "DENSVI_REPRES"
$ws.Range("E12").Value = "/ProcedimientoXunta/SI452A_2/ComprobacionDatos/tblDatos/FilaA2/cvDeniego"
$ws.Range("D12").Value = "Boolean"

# New row 13
$ws.Range("A13").Value = "autSviSol"
$ws.Range("B13").Value = "Autorización SVI solicitante"
$ws.Range("C13").Value = "Autorización SVI solicitante"
$ws.Range("F13").Value = "DENSVI_PRESENT"
$ws.Range("G13").Value = "DEN_SVI_SOL"
$ws.Range("E13").Value = "/ProcedimientoXunta/SI452A_2/ComprobacionDatos/tblDatos/FilaA2/cvDeniego"
$ws.Range("D13").Value = "Boolean"

# Fill G2:G11 with "no" for the existing rows (introduces "no" last)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Value = "no"
}

# Column widths (values chosen so the engine's internal pixel-quantized
# ColumnWidth lands on the target stored "width" attribute)
$ws.Columns.Item(5).ColumnWidth = 58.166666666666664
$ws.Columns.Item(7).ColumnWidth = 27.5

# Selection
[void]$ws.Range("G18").Select()
